$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.003.36'
$ws.Range('E2').Value = '  -2.35%  '
$ws.Range('D3').Value = '1.664.66'
$ws.Range('E3').Value = '  -1.91%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '216.50'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.35%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5091'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.72%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.005'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2630'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.99%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06385'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +2.12%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '21.78'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.72%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07402'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.07%  '
$ws.Range('D12').Value = '1.666.41'
$ws.Range('E12').Value = '  -1.89%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.496'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.37%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.5812'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.000008487'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.15'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -2.09%  '
$ws.Range('D17').Value = '26.057.84'
$ws.Range('E17').Value = '  -2.32%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '4.903'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.92%  '
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '10.66'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.02%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '188.88'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.85%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.194'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.22%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.006'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '145.28'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '7.573'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.55%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1187'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +3.23%  '
$ws.Range('E27').Value = '  -1.24%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.06672'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +16.70%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.301'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('E30').Value = '  -1.13%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.522'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.495'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.75%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.624'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.43%  '
$ws.Range('E34').Value = '  -0.60%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.6059'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.67%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.366'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.683'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.212'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +5.25%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01609'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.59%  '
$ws.Range('D40').Value = '1.075.46'
$ws.Range('E40').Value = '  -1.83%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.8586'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.76%  '
$ws.Range('E42').Value = '  +0.57%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '100.35'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.67%  '
$ws.Range('D44').Value = '1.811.57'
$ws.Range('E44').Value = '  -2.49%  '
$ws.Range('E45').Value = '  +3.71%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '56.20'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.09%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '8.017'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.78%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.05207'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.72%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.4290'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.80%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '5.940'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.97%  '
